$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; New="45÷5=9, 0"},
    @{Row=1;  Col=2; New="56÷5=11, 1"},
    @{Row=1;  Col=3; New="51÷8=6, 3"},
    @{Row=1;  Col=4; New="44÷4=11, 0"},
    @{Row=1;  Col=5; New="67÷3=22, 1"},

    @{Row=5;  Col=1; New="21÷2=10, 1"},
    @{Row=5;  Col=2; New="47÷5=9, 2"},
    @{Row=5;  Col=3; New="74÷9=8, 2"},
    @{Row=5;  Col=4; New="44÷6=7, 2"},
    @{Row=5;  Col=5; New="16÷4=4, 0"},

    @{Row=9;  Col=1; New="75÷6=12, 3"},
    @{Row=9;  Col=2; New="57÷3=19, 0"},
    @{Row=9;  Col=3; New="39÷6=6, 3"},
    @{Row=9;  Col=4; New="88÷9=9, 7"},
    @{Row=9;  Col=5; New="59÷7=8, 3"},

    @{Row=13; Col=1; New="53÷4=13, 1"},
    @{Row=13; Col=2; New="19÷5=3, 4"},
    @{Row=13; Col=3; New="43÷5=8, 3"},
    @{Row=13; Col=4; New="49÷8=6, 1"},
    @{Row=13; Col=5; New="14÷7=2, 0"},

    @{Row=17; Col=1; New="14÷7=2, 0"},
    @{Row=17; Col=2; New="64÷2=32, 0"},
    @{Row=17; Col=3; New="75÷6=12, 3"},
    @{Row=17; Col=4; New="27÷9=3, 0"},
    @{Row=17; Col=5; New="65÷5=13, 0"}
)

foreach ($u in $updates) {
    $cellRange = $t.Cell($u.Row, $u.Col).Range
    # Exclude the trailing end-of-cell marker so we don't destroy table structure
    $cellRange.End = $cellRange.End - 1
    $cellRange.Text = $u.New
}
